$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 19) - data shifted up by one period
$ws.Rows.Item(19).Delete()

# Rewrite rows 2-18 with the corrected naive-forecaster AR2 values

# Row 2
$ws.Cells.Item(2, 1).Value = 39765
$ws.Cells.Item(2, 2).Value = 2008
$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(2, 4).Value = 2009
$ws.Cells.Item(2, 5).ClearContents()

# Row 3
$ws.Cells.Item(3, 1).Value = 40130
$ws.Cells.Item(3, 2).Value = 2009
$ws.Cells.Item(3, 3).Value = 1.834695583582535
$ws.Cells.Item(3, 4).Value = 2010
$ws.Cells.Item(3, 5).ClearContents()

# Row 4
$ws.Cells.Item(4, 1).Value = 40494
$ws.Cells.Item(4, 2).Value = 2010
$ws.Cells.Item(4, 3).Value = 1.767835936772144
$ws.Cells.Item(4, 4).Value = 2011
$ws.Cells.Item(4, 5).ClearContents()

# Row 5
$ws.Cells.Item(5, 1).Value = 40862
$ws.Cells.Item(5, 2).Value = 2011
$ws.Cells.Item(5, 3).Value = 1.074400434091038
$ws.Cells.Item(5, 4).Value = 2012
$ws.Cells.Item(5, 5).ClearContents()

# Row 6
$ws.Cells.Item(6, 1).Value = 41228
$ws.Cells.Item(6, 2).Value = 2012
$ws.Cells.Item(6, 3).Value = 0.9212998022035679
$ws.Cells.Item(6, 4).Value = 2013
$ws.Cells.Item(6, 5).Value = 1.274704633957136

# Row 7
$ws.Cells.Item(7, 1).Value = 41592
$ws.Cells.Item(7, 2).Value = 2013
$ws.Cells.Item(7, 3).Value = 1.141837882844188
$ws.Cells.Item(7, 4).Value = 2014
$ws.Cells.Item(7, 5).Value = 1.404348988410131

# Row 8
$ws.Cells.Item(8, 1).Value = 41957
$ws.Cells.Item(8, 2).Value = 2014
$ws.Cells.Item(8, 3).Value = 1.335361538769475
$ws.Cells.Item(8, 4).Value = 2015
$ws.Cells.Item(8, 5).Value = 1.269653854937691

# Row 9
$ws.Cells.Item(9, 1).Value = 42321
$ws.Cells.Item(9, 2).Value = 2015
$ws.Cells.Item(9, 3).Value = 1.202048372526998
$ws.Cells.Item(9, 4).Value = 2016
$ws.Cells.Item(9, 5).Value = 1.253742200752095

# Row 10
$ws.Cells.Item(10, 1).Value = 42689
$ws.Cells.Item(10, 2).Value = 2016
$ws.Cells.Item(10, 3).Value = 2.677488680362305
$ws.Cells.Item(10, 4).Value = 2017
$ws.Cells.Item(10, 5).Value = 1.805615391969595

# Row 11
$ws.Cells.Item(11, 1).Value = 43053
$ws.Cells.Item(11, 2).Value = 2017
$ws.Cells.Item(11, 3).Value = 2.466954516646402
$ws.Cells.Item(11, 4).Value = 2018
$ws.Cells.Item(11, 5).Value = 1.661541796722577

# Row 12
$ws.Cells.Item(12, 1).Value = 43418
$ws.Cells.Item(12, 2).Value = 2018
$ws.Cells.Item(12, 3).Value = 1.401189216021326
$ws.Cells.Item(12, 4).Value = 2019
$ws.Cells.Item(12, 5).Value = 1.815016201748643

# Row 13
$ws.Cells.Item(13, 1).Value = 43783
$ws.Cells.Item(13, 2).Value = 2019
$ws.Cells.Item(13, 3).Value = 2.217567799050979
$ws.Cells.Item(13, 4).Value = 2020
$ws.Cells.Item(13, 5).Value = 1.810449264563152

# Row 14
$ws.Cells.Item(14, 1).Value = 44159
$ws.Cells.Item(14, 2).Value = 2020
$ws.Cells.Item(14, 3).Value = 2.139672475020404
$ws.Cells.Item(14, 4).Value = 2021
$ws.Cells.Item(14, 5).Value = 2.128328071999674

# Row 15
$ws.Cells.Item(15, 1).Value = 44525
$ws.Cells.Item(15, 2).Value = 2021
$ws.Cells.Item(15, 3).Value = 2.100991693542231
$ws.Cells.Item(15, 4).Value = 2022
$ws.Cells.Item(15, 5).Value = 1.11435041103376

# Row 16
$ws.Cells.Item(16, 1).Value = 44890
$ws.Cells.Item(16, 2).Value = 2022
$ws.Cells.Item(16, 3).Value = 0.8967077601845341
$ws.Cells.Item(16, 4).Value = 2023
$ws.Cells.Item(16, 5).Value = 0.3338851812143995

# Row 17
$ws.Cells.Item(17, 1).Value = 45254
$ws.Cells.Item(17, 2).Value = 2023
$ws.Cells.Item(17, 3).Value = 0.782207885866093
$ws.Cells.Item(17, 4).Value = 2024
$ws.Cells.Item(17, 5).Value = 2.228542839642689

# Row 18
$ws.Cells.Item(18, 1).Value = 45618
$ws.Cells.Item(18, 2).Value = 2024
$ws.Cells.Item(18, 3).Value = 1.508385007449875
$ws.Cells.Item(18, 4).Value = 2025
$ws.Cells.Item(18, 5).Value = 0.9823016603409229
